# Auto-generated script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.201.83"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.057.38"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.26"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  +3.25%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.14"
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.88"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0759"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.358.48"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.56"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.82"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.778"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.056.99"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.141.34"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.41"
$ws.Range("E20").Value = "  +9.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.22"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0809"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.77"
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.44"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("E28").Value = "  +7.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.78"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.45"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0617"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("E35").Value = "  +7.48%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.74"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.26"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.68"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.40"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.471.14"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.32"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("E45").Value = "  +5.61%  "
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.09"
$ws.Range("E49").Value = "  -3.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("E51").Value = "  +1.39%  "
